$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking values must stay text (match source formatting,
# e.g. keep trailing zeros / avoid Excel auto-converting "0.9998" to a number).
$textCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.294.30'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.865.50'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").Value = '311.05'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").Value = '0.9989'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").Value = '0.4986'
$ws.Range("E7").Value = '  -3.26%  '
$ws.Range("D8").Value = '0.3981'
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("D9").Value = '0.1014'
$ws.Range("E9").Value = '  +29.66%  '
$ws.Range("D10").Value = '1.122'
$ws.Range("E10").Value = '  +1.00%  '
$ws.Range("D11").Value = '41.46'
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").Value = '6.494'
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("D13").Value = '21.06'
$ws.Range("E13").Value = '  +3.03%  '
$ws.Range("D14").Value = '1.869.52'
$ws.Range("E14").Value = '  +3.52%  '
$ws.Range("D15").Value = '0.9998'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '7.385'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").Value = '0.00001144'
$ws.Range("E17").Value = '  +5.69%  '
$ws.Range("D18").Value = '93.69'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").Value = '0.06653'
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").Value = '17.37'
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = '6.091'
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("D23").Value = '28.368.83'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").Value = '11.36'
$ws.Range("E24").Value = '  +2.03%  '
$ws.Range("D25").Value = '2.239'
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.508'
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '21.19'
$ws.Range("E27").Value = '  +3.32%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.076.98'
$ws.Range("E28").Value = '  +2.92%  '
$ws.Range("D29").Value = '157.67'
$ws.Range("E29").Value = '  -1.90%  '
$ws.Range("D30").Value = '127.75'
$ws.Range("E30").Value = '  -0.73%  '
$ws.Range("D31").Value = '0.1060'
$ws.Range("E31").Value = '  -3.80%  '
$ws.Range("D32").Value = '1.059'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '5.640'
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").Value = '3.595'
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").Value = '0.06808'
$ws.Range("E35").Value = '  -5.21%  '
$ws.Range("D36").Value = '9.249'
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("D37").Value = '0.02399'
$ws.Range("E37").Value = '  +1.98%  '
$ws.Range("D38").Value = '0.2183'
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("D39").Value = '5.023'
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").Value = '0.6302'
$ws.Range("E41").Value = '  +2.01%  '
$ws.Range("D42").Value = '1.177'
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("D43").Value = '0.9998'
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = '13.40'
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").Value = '0.5998'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '1.280'
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("D47").Value = '3.675'
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("D48").Value = '125.04'
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("D49").Value = '1.992'
$ws.Range("E49").Value = '  +3.73%  '
$ws.Range("D50").Value = '1.192'
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("D51").Value = '1.120'
$ws.Range("E51").Value = '  +4.45%  '
